$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"0.90715"
$ws.Cells.Item(2, 8).Value = [double]"2.72145"
$ws.Cells.Item(2, 9).Value = [double]"0.01717809939998381"
$ws.Cells.Item(2, 10).Value = [double]"0.01717809939998381"
$ws.Cells.Item(2, 13).Value = [double]"1.461859"
$ws.Cells.Item(2, 14).Value = [double]"4.385577"
$ws.Cells.Item(2, 15).Value = [double]"0.1790970628189019"
$ws.Cells.Item(2, 16).Value = [double]"0.1790970628189019"
$ws.Cells.Item(2, 17).Value = [double]"1.32612539185"
$ws.Cells.Item(2, 18).Value = [double]"11.93512852665"
$ws.Cells.Item(2, 19).Value = [double]"0.003076547147348242"
$ws.Cells.Item(2, 20).Value = [double]"0.003076547147348241"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"0.90715"
$ws.Cells.Item(3, 8).Value = [double]"2.72145"
$ws.Cells.Item(3, 9).Value = [double]"0.01717809939998381"
$ws.Cells.Item(3, 10).Value = [double]"0.01717809939998381"
$ws.Cells.Item(3, 15).Value = [double]"0.1506444277390854"
$ws.Cells.Item(3, 16).Value = [double]"0.1506444277390854"
$ws.Cells.Item(3, 17).Value = [double]"1.115447666316667"
$ws.Cells.Item(3, 18).Value = [double]"10.03902899685"
$ws.Cells.Item(3, 19).Value = [double]"0.002587784953755688"
$ws.Cells.Item(3, 20).Value = [double]"0.002587784953755686"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"0.90715"
$ws.Cells.Item(4, 8).Value = [double]"2.72145"
$ws.Cells.Item(4, 9).Value = [double]"0.01717809939998381"
$ws.Cells.Item(4, 10).Value = [double]"0.01717809939998381"
$ws.Cells.Item(4, 13).Value = [double]"5.443148666666667"
$ws.Cells.Item(4, 14).Value = [double]"16.329446"
$ws.Cells.Item(4, 15).Value = [double]"0.6668577056245659"
$ws.Cells.Item(4, 16).Value = [double]"0.6668577056245659"
$ws.Cells.Item(4, 17).Value = [double]"4.937752312966667"
$ws.Cells.Item(4, 18).Value = [double]"44.4397708167"
$ws.Cells.Item(4, 19).Value = [double]"0.01145534795286394"
$ws.Cells.Item(4, 20).Value = [double]"0.01145534795286393"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"0.90715"
$ws.Cells.Item(5, 8).Value = [double]"2.72145"
$ws.Cells.Item(5, 9).Value = [double]"0.01717809939998381"
$ws.Cells.Item(5, 10).Value = [double]"0.01717809939998381"
$ws.Cells.Item(5, 13).Value = [double]"0.02775866666666667"
$ws.Cells.Item(5, 14).Value = [double]"0.083276"
$ws.Cells.Item(5, 15).Value = [double]"0.003400803817446799"
$ws.Cells.Item(5, 16).Value = [double]"0.003400803817446798"
$ws.Cells.Item(5, 17).Value = [double]"0.02518127446666667"
$ws.Cells.Item(5, 18).Value = [double]"0.2266314702"
$ws.Cells.Item(5, 19).Value = [double]"5.841934601594551E-05"
$ws.Cells.Item(5, 20).Value = [double]"5.841934601594549E-05"
$ws.Cells.Item(6, 9).Value = [double]"0.0237180037344858"
$ws.Cells.Item(6, 10).Value = [double]"0.0237180037344858"
$ws.Cells.Item(6, 13).Value = [double]"1.461859"
$ws.Cells.Item(6, 14).Value = [double]"4.385577"
$ws.Cells.Item(6, 15).Value = [double]"0.1790970628189019"
$ws.Cells.Item(6, 16).Value = [double]"0.1790970628189019"
$ws.Cells.Item(6, 17).Value = [double]"1.830996914380667"
$ws.Cells.Item(6, 18).Value = [double]"16.478972229426"
$ws.Cells.Item(6, 19).Value = [double]"0.004247824804774153"
$ws.Cells.Item(6, 20).Value = [double]"0.004247824804774153"
$ws.Cells.Item(7, 9).Value = [double]"0.0237180037344858"
$ws.Cells.Item(7, 10).Value = [double]"0.0237180037344858"
$ws.Cells.Item(7, 15).Value = [double]"0.1506444277390854"
$ws.Cells.Item(7, 16).Value = [double]"0.1506444277390854"
$ws.Cells.Item(7, 19).Value = [double]"0.003572985099695104"
$ws.Cells.Item(7, 20).Value = [double]"0.003572985099695103"
$ws.Cells.Item(8, 9).Value = [double]"0.0237180037344858"
$ws.Cells.Item(8, 10).Value = [double]"0.0237180037344858"
$ws.Cells.Item(8, 13).Value = [double]"5.443148666666667"
$ws.Cells.Item(8, 14).Value = [double]"16.329446"
$ws.Cells.Item(8, 15).Value = [double]"0.6668577056245659"
$ws.Cells.Item(8, 16).Value = [double]"0.6668577056245659"
$ws.Cells.Item(8, 17).Value = [double]"6.817612651549779"
$ws.Cells.Item(8, 18).Value = [double]"61.358513863948"
$ws.Cells.Item(8, 19).Value = [double]"0.01581653355237409"
$ws.Cells.Item(8, 20).Value = [double]"0.01581653355237408"
$ws.Cells.Item(9, 9).Value = [double]"0.0237180037344858"
$ws.Cells.Item(9, 10).Value = [double]"0.0237180037344858"
$ws.Cells.Item(9, 13).Value = [double]"0.02775866666666667"
$ws.Cells.Item(9, 14).Value = [double]"0.083276"
$ws.Cells.Item(9, 15).Value = [double]"0.003400803817446799"
$ws.Cells.Item(9, 16).Value = [double]"0.003400803817446798"
$ws.Cells.Item(9, 17).Value = [double]"0.03476808160977778"
$ws.Cells.Item(9, 18).Value = [double]"0.3129127344880001"
$ws.Cells.Item(9, 19).Value = [double]"8.066027764245673E-05"
$ws.Cells.Item(9, 20).Value = [double]"8.066027764245672E-05"
$ws.Cells.Item(10, 7).Value = [double]"2.247832333333333"
$ws.Cells.Item(10, 8).Value = [double]"6.743497"
$ws.Cells.Item(10, 9).Value = [double]"0.04256571378106988"
$ws.Cells.Item(10, 10).Value = [double]"0.04256571378106987"
$ws.Cells.Item(10, 13).Value = [double]"1.461859"
$ws.Cells.Item(10, 14).Value = [double]"4.385577"
$ws.Cells.Item(10, 15).Value = [double]"0.1790970628189019"
$ws.Cells.Item(10, 16).Value = [double]"0.1790970628189019"
$ws.Cells.Item(10, 17).Value = [double]"3.286013926974332"
$ws.Cells.Item(10, 18).Value = [double]"29.57412534276899"
$ws.Cells.Item(10, 19).Value = [double]"0.007623394314979671"
$ws.Cells.Item(10, 20).Value = [double]"0.007623394314979669"
$ws.Cells.Item(11, 7).Value = [double]"2.247832333333333"
$ws.Cells.Item(11, 8).Value = [double]"6.743497"
$ws.Cells.Item(11, 9).Value = [double]"0.04256571378106988"
$ws.Cells.Item(11, 10).Value = [double]"0.04256571378106987"
$ws.Cells.Item(11, 15).Value = [double]"0.1506444277390854"
$ws.Cells.Item(11, 16).Value = [double]"0.1506444277390854"
$ws.Cells.Item(11, 17).Value = [double]"2.763974348771222"
$ws.Cells.Item(11, 18).Value = [double]"24.875769138941"
$ws.Cells.Item(11, 19).Value = [double]"0.006412287593854973"
$ws.Cells.Item(11, 20).Value = [double]"0.006412287593854971"
$ws.Cells.Item(12, 7).Value = [double]"2.247832333333333"
$ws.Cells.Item(12, 8).Value = [double]"6.743497"
$ws.Cells.Item(12, 9).Value = [double]"0.04256571378106988"
$ws.Cells.Item(12, 10).Value = [double]"0.04256571378106987"
$ws.Cells.Item(12, 13).Value = [double]"5.443148666666667"
$ws.Cells.Item(12, 14).Value = [double]"16.329446"
$ws.Cells.Item(12, 15).Value = [double]"0.6668577056245659"
$ws.Cells.Item(12, 16).Value = [double]"0.6668577056245659"
$ws.Cells.Item(12, 17).Value = [double]"12.23528556807356"
$ws.Cells.Item(12, 18).Value = [double]"110.117570112662"
$ws.Cells.Item(12, 19).Value = [double]"0.02838527423031623"
$ws.Cells.Item(12, 20).Value = [double]"0.02838527423031622"
$ws.Cells.Item(13, 7).Value = [double]"2.247832333333333"
$ws.Cells.Item(13, 8).Value = [double]"6.743497"
$ws.Cells.Item(13, 9).Value = [double]"0.04256571378106988"
$ws.Cells.Item(13, 10).Value = [double]"0.04256571378106987"
$ws.Cells.Item(13, 13).Value = [double]"0.02775866666666667"
$ws.Cells.Item(13, 14).Value = [double]"0.083276"
$ws.Cells.Item(13, 15).Value = [double]"0.003400803817446799"
$ws.Cells.Item(13, 16).Value = [double]"0.003400803817446798"
$ws.Cells.Item(13, 17).Value = [double]"0.06239682846355556"
$ws.Cells.Item(13, 18).Value = [double]"0.561571456172"
$ws.Cells.Item(13, 19).Value = [double]"0.0001447576419190103"
$ws.Cells.Item(13, 20).Value = [double]"0.0001447576419190102"
$ws.Cells.Item(14, 7).Value = [double]"48.40102466666667"
$ws.Cells.Item(14, 8).Value = [double]"145.203074"
$ws.Cells.Item(14, 9).Value = [double]"0.9165381830844606"
$ws.Cells.Item(14, 10).Value = [double]"0.9165381830844604"
$ws.Cells.Item(14, 13).Value = [double]"1.461859"
$ws.Cells.Item(14, 14).Value = [double]"4.385577"
$ws.Cells.Item(14, 15).Value = [double]"0.1790970628189019"
$ws.Cells.Item(14, 16).Value = [double]"0.1790970628189019"
$ws.Cells.Item(14, 17).Value = [double]"70.75547351818867"
$ws.Cells.Item(14, 18).Value = [double]"636.799261663698"
$ws.Cells.Item(14, 19).Value = [double]"0.1641492965517999"
$ws.Cells.Item(14, 20).Value = [double]"0.1641492965517998"
$ws.Cells.Item(15, 7).Value = [double]"48.40102466666667"
$ws.Cells.Item(15, 8).Value = [double]"145.203074"
$ws.Cells.Item(15, 9).Value = [double]"0.9165381830844606"
$ws.Cells.Item(15, 10).Value = [double]"0.9165381830844604"
$ws.Cells.Item(15, 15).Value = [double]"0.1506444277390854"
$ws.Cells.Item(15, 16).Value = [double]"0.1506444277390854"
$ws.Cells.Item(15, 17).Value = [double]"59.51475501490246"
$ws.Cells.Item(15, 18).Value = [double]"535.632795134122"
$ws.Cells.Item(15, 19).Value = [double]"0.1380713700917796"
$ws.Cells.Item(15, 20).Value = [double]"0.1380713700917796"
$ws.Cells.Item(16, 7).Value = [double]"48.40102466666667"
$ws.Cells.Item(16, 8).Value = [double]"145.203074"
$ws.Cells.Item(16, 9).Value = [double]"0.9165381830844606"
$ws.Cells.Item(16, 10).Value = [double]"0.9165381830844604"
$ws.Cells.Item(16, 13).Value = [double]"5.443148666666667"
$ws.Cells.Item(16, 14).Value = [double]"16.329446"
$ws.Cells.Item(16, 15).Value = [double]"0.6668577056245659"
$ws.Cells.Item(16, 16).Value = [double]"0.6668577056245659"
$ws.Cells.Item(16, 17).Value = [double]"263.4539728796672"
$ws.Cells.Item(16, 18).Value = [double]"2371.085755917004"
$ws.Cells.Item(16, 19).Value = [double]"0.6112005498890117"
$ws.Cells.Item(16, 20).Value = [double]"0.6112005498890116"
$ws.Cells.Item(17, 7).Value = [double]"48.40102466666667"
$ws.Cells.Item(17, 8).Value = [double]"145.203074"
$ws.Cells.Item(17, 9).Value = [double]"0.9165381830844606"
$ws.Cells.Item(17, 10).Value = [double]"0.9165381830844604"
$ws.Cells.Item(17, 13).Value = [double]"0.02775866666666667"
$ws.Cells.Item(17, 14).Value = [double]"0.083276"
$ws.Cells.Item(17, 15).Value = [double]"0.003400803817446799"
$ws.Cells.Item(17, 16).Value = [double]"0.003400803817446798"
$ws.Cells.Item(17, 17).Value = [double]"1.343547910047111"
$ws.Cells.Item(17, 18).Value = [double]"12.091931190424"
$ws.Cells.Item(17, 19).Value = [double]"0.003116966551869386"
$ws.Cells.Item(17, 20).Value = [double]"0.003116966551869385"
